# Apply the changes described by the diff:
# 1. Update a few odds values on row 2 and row 4
# 2. Delete rows 5 and 6 entirely (their matches were removed)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10

# Row 4 changes
$ws.Range("I4").Value = 2.18
$ws.Range("K4").Value = 2.15
$ws.Range("L4").Value = 2.72
$ws.Range("Y4").Value = 10.5
$ws.Range("Z4").Value = 37
$ws.Range("AD4").Value = 6.7
$ws.Range("AE4").Value = 12
$ws.Range("AK4").Value = 22
$ws.Range("AN4").Value = 5
$ws.Range("AR4").Value = 90
$ws.Range("AX4").Value = 11

# Delete rows 5 and 6 (shift rows up)
$ws.Range("A5:BD6").EntireRow.Delete()
